# TASK_42, completed, rating realy loading right now, mmm...
#
# - sheet1 ("Задачи"): mark completion date for row 45 (E45) and start
#   date for row 46 (D46), matching the style already used by the other
#   date cells in that block.
# - sheet2 ("Бэклог задач"): append four new backlog items discovered
#   while finishing the rating work.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: "Задачи" -------------------------------------------------

# Row 45 ("10.3 – Обновление рейтинга на клиенте.") finished.
$ws1.Range("E45").Value = 42020.774305555555

# Row 46 ("10.4 – Расчёт рейтинга.") started. D46 had no date format
# applied yet, so copy the formatting from the cell above (D45) before
# writing the value, same as Excel would do when you fill this in.
$ws1.Range("D45").Copy() | Out-Null
$ws1.Range("D46").PasteSpecial(-4122) | Out-Null
$ws1.Range("D46").Value = 42020.774305555555

$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 38
$ws1.Range("F38").Select() | Out-Null

# --- Sheet2: "Бэклог задач" --------------------------------------------

# New backlog rows, following on from the rating work. Copy the number
# format from the last existing date cell (C37) onto the new date cells
# first, so they pick up the same style.
$ws2.Range("C37").Copy() | Out-Null
$ws2.Range("C38:C41").PasteSpecial(-4122) | Out-Null

$ws2.Range("B38").Value = "Оптимизировать загрузку рейтинга. Сейчас при заходе на страницу рейтинга, обновляется каждую секунду, в т.ч. Данные очков берутьс из юзера, а вот позиции могут не соответстовать."
$ws2.Range("C38").Value = 42020.765972222223

$ws2.Range("B39").Value = "Сделать заголовок рейтингу. Его нет, а это печально, прям жесть, без этого нельзя запускать."
$ws2.Range("C39").Value = 42020.76666666667

$ws2.Range("B40").Value = "Добавить рейтинг по друзьям."
$ws2.Range("C40").Value = 42020.767361111109

$ws2.Range("B41").Value = "Добавить рейтинг по позиции игрока."
$ws2.Range("C41").Value = 42020.767361111109

$ws2.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$ws2.Range("B41").Select() | Out-Null

# Leave sheet1 as the tab that is active/selected when the file is saved,
# matching the original workbook (tabSelected="1" on sheet1).
$ws1.Activate() | Out-Null
